$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert two fresh rows directly below row 2 (they become rows 3 & 4) ---
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# --- Step 2: re-create the original record (previously row 2) in row 3, unchanged ---
$ws.Range("A3").Value = 114565031
$ws.Range("B3").Value = 78507
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = "Garnlav"
$ws.Range("G3").Value = "Alectoria sarmentosa"
$ws.Range("H3").Value = "(Ach.) Ach."
$ws.Range("P3").Value = "Horndal-Fallviken, Dlr"
$ws.Range("Q3").Value = 565454
$ws.Range("R3").Value = 6746042
$ws.Range("S3").Value = 10
$ws.Range("T3").Value = "Dalarna"
$ws.Range("U3").Value = "Falun"
$ws.Range("V3").Value = "Dalarna"
$ws.Range("W3").Value = "Svärdsjö"
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2023-05-10"
$ws.Range("Z3").Value = "13:14"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2023-05-10"
$ws.Range("AB3").Value = "13:14"
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AW3").Value = "Mårten Nilsson"
$ws.Range("AX3").Value = "Ward Tamsyn"

# --- Step 3: update row 2 in place with the corrected/updated record (Tallticka) ---
$ws.Range("A2").Value = 108951146
$ws.Range("B2").Value = 89412
$ws.Range("E2").Value = 5442
$ws.Range("F2").Value = "Tallticka"
$ws.Range("G2").Value = "Porodaedalea pini"
$ws.Range("H2").Value = "(Brot.) Murrill"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "1"
$ws.Range("J2").Value = ""
$ws.Range("N2").Value = ""
$ws.Range("P2").Value = "Stora Mörtsjön, Dlr"
$ws.Range("Q2").Value = 565429.1409946628
$ws.Range("R2").Value = 6746061.659947474
$ws.Range("S2").Value = 25
$ws.Range("Z2").Value = "15:12"
$ws.Range("AB2").Value = "15:12"
$ws.Range("AF2").Value = ""
$ws.Range("AJ2").Value = "tall"
$ws.Range("AK2").Value = "Pinus sylvestris"
$ws.Range("AO2").Value = "Pinus sylvestris"
$ws.Range("AW2").Value = "Ward Tamsyn"

# --- Step 4: populate row 4 with the brand new second "Tallticka" record ---
$ws.Range("A4").Value = 114571384
$ws.Range("B4").Value = 90352
$ws.Range("C4").Value = "Ovaliderad"
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 5442
$ws.Range("F4").Value = "Tallticka"
$ws.Range("G4").Value = "Porodaedalea pini"
$ws.Range("H4").Value = "(Brot.) Murrill"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "1"
$ws.Range("P4").Value = "Horndal-Fallviken, Dlr"
$ws.Range("Q4").Value = 565455
$ws.Range("R4").Value = 6746088
$ws.Range("S4").Value = 10
$ws.Range("T4").Value = "Dalarna"
$ws.Range("U4").Value = "Falun"
$ws.Range("V4").Value = "Dalarna"
$ws.Range("W4").Value = "Svärdsjö"
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2023-05-10"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2023-05-10"
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AW4").Value = "Mårten Nilsson"
$ws.Range("AX4").Value = "Ward Tamsyn"
